# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Chocobo_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the target commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 700.5
$ws.Range("I2").Value = 700.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 700.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -587.5
$ws.Range("H18").Value = 435
$ws.Range("I18").Value = 325
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 325
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = -41
$ws.Range("N18").Value = -1168
$ws.Range("H96").Value = 612.6667
$ws.Range("I96").Value = 299.23077
$ws.Range("K96").Value = 897.69231
$ws.Range("M96").Value = 475.30769
$ws.Range("H100").Value = 40001880
$ws.Range("I100").Value = 66668300
$ws.Range("J100").Value = 2250
$ws.Range("K100").Value = 66668300
$ws.Range("L100").Value = 2250
$ws.Range("M100").Value = -66667759
$ws.Range("N100").Value = -3332
$ws.Range("H113").Value = 10520.75
$ws.Range("I113").Value = 4057
$ws.Range("K113").Value = 4057
$ws.Range("M113").Value = -803
$ws.Range("H116").Value = 1256812.5
$ws.Range("I116").Value = 2501500
$ws.Range("J116").Value = 12125
$ws.Range("K116").Value = 2501500
$ws.Range("L116").Value = 12125
$ws.Range("M116").Value = -2498058
$ws.Range("N116").Value = -19009
$ws.Range("H132").Value = 35860064
$ws.Range("I132").Value = 40162190
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 120486570
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -120484040
$ws.Range("N132").Value = -32060
$ws.Range("H137").Value = 3877.3096
$ws.Range("I137").Value = 3182.9429
$ws.Range("K137").Value = 9548.8287
$ws.Range("M137").Value = -6998.8287
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 7323.778
$ws.Range("I3").Value = 2802
$ws.Range("J3").Value = 12976
$ws.Range("K3").Value = 2802
$ws.Range("L3").Value = 12976
$ws.Range("M3").Value = -2687
$ws.Range("N3").Value = -13206
$ws.Range("H25").Value = 4318.7144
$ws.Range("I25").Value = 1646.4
$ws.Range("J25").Value = 10999.5
$ws.Range("K25").Value = 1646.4
$ws.Range("L25").Value = 10999.5
$ws.Range("M25").Value = -1244.4
$ws.Range("N25").Value = -11803.5
$ws.Range("H32").Value = 10047.389
$ws.Range("I32").Value = 7110.0566
$ws.Range("J32").Value = 14912.344
$ws.Range("K32").Value = 7110.0566
$ws.Range("L32").Value = 14912.344
$ws.Range("M32").Value = -6823.0566
$ws.Range("N32").Value = -15486.344
$ws.Range("H74").Value = 2313.8572
$ws.Range("I74").Value = 1693.2307
$ws.Range("J74").Value = 3322.375
$ws.Range("K74").Value = 1693.2307
$ws.Range("L74").Value = 3322.375
$ws.Range("M74").Value = -819.2307000000001
$ws.Range("N74").Value = -5070.375
$ws.Range("H76").Value = 26686
$ws.Range("J76").Value = 26686
$ws.Range("L76").Value = 26686
$ws.Range("N76").Value = -27362
$ws.Range("H77").Value = 2313.8572
$ws.Range("I77").Value = 1693.2307
$ws.Range("J77").Value = 3322.375
$ws.Range("K77").Value = 8466.1535
$ws.Range("L77").Value = 16611.875
$ws.Range("M77").Value = -4098.1535
$ws.Range("N77").Value = -25347.875
$ws.Range("H79").Value = 26686
$ws.Range("J79").Value = 26686
$ws.Range("L79").Value = 26686
$ws.Range("N79").Value = -29026
$ws.Range("H122").Value = 2654
$ws.Range("I122").Value = 1483.5555
$ws.Range("J122").Value = 5287.5
$ws.Range("K122").Value = 4450.666499999999
$ws.Range("L122").Value = 15862.5
$ws.Range("M122").Value = -2000.666499999999
$ws.Range("N122").Value = -20762.5
$ws.Range("H132").Value = 2600.5806
$ws.Range("I132").Value = 1356.3182
$ws.Range("K132").Value = 4068.9546
$ws.Range("M132").Value = -1538.9546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1307.0714
$ws.Range("I5").Value = 1030
$ws.Range("J5").Value = 1584.1428
$ws.Range("K5").Value = 1030
$ws.Range("L5").Value = 1584.1428
$ws.Range("M5").Value = -917
$ws.Range("N5").Value = -1810.1428
$ws.Range("H82").Value = 16623.45
$ws.Range("I82").Value = 2418.6667
$ws.Range("J82").Value = 28245.545
$ws.Range("K82").Value = 2418.6667
$ws.Range("L82").Value = 28245.545
$ws.Range("M82").Value = -2035.6667
$ws.Range("N82").Value = -29011.545
$ws.Range("H85").Value = 16623.45
$ws.Range("I85").Value = 2418.6667
$ws.Range("J85").Value = 28245.545
$ws.Range("K85").Value = 2418.6667
$ws.Range("L85").Value = 28245.545
$ws.Range("M85").Value = -1092.6667
$ws.Range("N85").Value = -30897.545
$ws.Range("H134").Value = 3925.5
$ws.Range("I134").Value = 2116.762
$ws.Range("J134").Value = 8145.8887
$ws.Range("K134").Value = 6350.286
$ws.Range("L134").Value = 24437.6661
$ws.Range("M134").Value = -3815.286
$ws.Range("N134").Value = -29507.6661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 5149.5
$ws.Range("I2").Value = 5149.5
$ws.Range("K2").Value = 5149.5
$ws.Range("M2").Value = -5036.5
$ws.Range("H31").Value = 3595.366
$ws.Range("I31").Value = 1369.4348
$ws.Range("J31").Value = 6439.6113
$ws.Range("K31").Value = 1369.4348
$ws.Range("L31").Value = 6439.6113
$ws.Range("M31").Value = -1074.4348
$ws.Range("N31").Value = -7029.6113
$ws.Range("H34").Value = 3595.366
$ws.Range("I34").Value = 1369.4348
$ws.Range("J34").Value = 6439.6113
$ws.Range("K34").Value = 1369.4348
$ws.Range("L34").Value = 6439.6113
$ws.Range("M34").Value = -1167.4348
$ws.Range("N34").Value = -6843.6113
$ws.Range("H35").Value = 1913
$ws.Range("I35").Value = 1913
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1913
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1619
$ws.Range("H68").Value = 50167.668
$ws.Range("J68").Value = 50167.668
$ws.Range("L68").Value = 50167.668
$ws.Range("N68").Value = -51665.668
$ws.Range("H71").Value = 50167.668
$ws.Range("J71").Value = 50167.668
$ws.Range("L71").Value = 150503.004
$ws.Range("N71").Value = -157991.004
$ws.Range("H132").Value = 3136.2122
$ws.Range("I132").Value = 2536.3215
$ws.Range("J132").Value = 6495.6
$ws.Range("K132").Value = 7608.9645
$ws.Range("L132").Value = 19486.8
$ws.Range("M132").Value = -5078.9645
$ws.Range("N132").Value = -24546.8
$ws.Range("H141").Value = 18141.463
$ws.Range("J141").Value = 18141.463
$ws.Range("L141").Value = 18141.463
$ws.Range("N141").Value = -28501.463
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6027643.5
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 3000
$ws.Range("N4").Value = -3224
$ws.Range("H113").Value = 749.2683
$ws.Range("I113").Value = 645.73334
$ws.Range("K113").Value = 1937.20002
$ws.Range("M113").Value = 232.79998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 29124.875
$ws.Range("J4").Value = 29124.875
$ws.Range("L4").Value = 29124.875
$ws.Range("N4").Value = -29348.875
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H102").Value = 2889.4075
$ws.Range("I102").Value = 2338.0833
$ws.Range("J102").Value = 7300
$ws.Range("K102").Value = 2338.0833
$ws.Range("L102").Value = 7300
$ws.Range("M102").Value = -716.0832999999998
$ws.Range("N102").Value = -10544
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10002127
$ws.Range("I22").Value = 16668201
$ws.Range("J22").Value = 3015.1
$ws.Range("K22").Value = 16668201
$ws.Range("L22").Value = 3015.1
$ws.Range("M22").Value = -16667906
$ws.Range("N22").Value = -3605.1
$ws.Range("H27").Value = 10002127
$ws.Range("I27").Value = 16668201
$ws.Range("J27").Value = 3015.1
$ws.Range("K27").Value = 16668201
$ws.Range("L27").Value = 3015.1
$ws.Range("M27").Value = -16668094
$ws.Range("N27").Value = -3229.1
$ws.Range("H40").Value = 6146.1724
$ws.Range("I40").Value = 5964.7144
$ws.Range("J40").Value = 6622.5
$ws.Range("K40").Value = 5964.7144
$ws.Range("L40").Value = 6622.5
$ws.Range("M40").Value = -5828.7144
$ws.Range("N40").Value = -6894.5
$ws.Range("H74").Value = 32998.918
$ws.Range("J74").Value = 40309.89
$ws.Range("L74").Value = 40309.89
$ws.Range("N74").Value = -42305.89
$ws.Range("H77").Value = 32998.918
$ws.Range("J77").Value = 40309.89
$ws.Range("L77").Value = 120929.67
$ws.Range("N77").Value = -130913.67
$ws.Range("H122").Value = 6422.1665
$ws.Range("I122").Value = 5128.4287
$ws.Range("J122").Value = 7245.4546
$ws.Range("K122").Value = 15385.2861
$ws.Range("L122").Value = 21736.3638
$ws.Range("M122").Value = -12935.2861
$ws.Range("N122").Value = -26636.3638
$ws.Range("H132").Value = 6611.355
$ws.Range("I132").Value = 2456.2856
$ws.Range("J132").Value = 7823.25
$ws.Range("K132").Value = 7368.8568
$ws.Range("L132").Value = 23469.75
$ws.Range("M132").Value = -4838.8568
$ws.Range("N132").Value = -28529.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11000
$ws.Range("I122").Value = 7000
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 21000
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -18550
$ws.Range("N122").Value = -49900
$ws.Range("H132").Value = 8337715
$ws.Range("I132").Value = 5399.1904
$ws.Range("J132").Value = 17547118
$ws.Range("K132").Value = 16197.5712
$ws.Range("L132").Value = 52641354
$ws.Range("M132").Value = -13667.5712
$ws.Range("N132").Value = -52646414
$ws.Range("H136").Value = 3164.25
$ws.Range("I136").Value = 887.9091
$ws.Range("K136").Value = 2663.7273
$ws.Range("M136").Value = -113.7273
